# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Refreshes the Micro/SMEs/MSMEs percentage figures on the Estonia Summary
# sheet with more precise (two decimal place) values, for both the
# "Statistical Institution" source table and the "SME Associations" source
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be stored as Text so the precision-bearing values
    # (e.g. "44.33") are not silently reinterpreted as numbers.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# "Source Type: Statistical Institution" table
# Row 11: Enterprises density (per 1000 people)
Set-TextValue "B11" "44.33"
Set-TextValue "C11" "4.88"
Set-TextValue "D11" "49.21"

# Row 12: Employment (% of total)
Set-TextValue "B12" "27.58"
Set-TextValue "C12" "50.84"
Set-TextValue "D12" "78.41"

# "Source Type: SME Associations (Most Widely Used)" table
# Row 33: Enterprises density (per 1000 people)
Set-TextValue "B33" "36.63"
Set-TextValue "C33" "4.26"
Set-TextValue "D33" "40.89"

# Row 34: Employment (% of total)
Set-TextValue "B34" "29.06"
Set-TextValue "C34" "48.99"
Set-TextValue "D34" "78.05"

# Row 36: Enterprises (% of total)
Set-TextValue "B36" "89.34"
Set-TextValue "D36" "99.73"
